$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "榮" character that sits between "世尊" and
#    "藥師" in the first paragraph's heading line, turning
#    "頂禮世尊榮藥師琉璃光王如來應正等覺" into
#    "頂禮世尊藥師琉璃光王如來應正等覺".
# ------------------------------------------------------------------
$findRange = $d.Content
[void]$findRange.Find.Execute("榮", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rongStart = $findRange.Start
$rongRange = $d.Range($rongStart, $rongStart + 1)
$rongRange.Delete()

# ------------------------------------------------------------------
# 2. Split the now-merged run into two runs - "頂禮世尊" and
#    "藥師琉璃光王如來應正等覺" - that carry identical character
#    formatting. Word only materializes a new run boundary when a
#    sub-range's formatting is actually touched, so toggle Bold on
#    the trailing half and immediately revert it.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$secondHalf = $d.Range($rongStart, $p1.Range.End - 1)
$secondHalf.Font.Bold = 1
$secondHalf.Font.Bold = 0

# ------------------------------------------------------------------
# 3. Relocate the hidden "_GoBack" bookmark. It currently wraps the
#    whole mantra paragraph ("tädyātha ...svāhā"); move it so it
#    becomes an empty (collapsed) bookmark at the very start of the
#    blank paragraph that precedes the mantra.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$blankPara = $d.Paragraphs(2).Range
$blankPara.Collapse(1)
$d.Bookmarks.Add("_GoBack", $blankPara)
